$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23 (邱志鹏): plan content + completion status
$ws.Range("B23").Value = "初步设计完成APP端用户登陆、群聊创建的用例分析和界面设计          "
$ws.Range("C23").Value = "未完成（90%）"

# Row 24 (黄立根): plan content + completion status
$ws.Range("B24").Value = "初步设计完成网页端的登陆、注册、主页面、个人信息管理的界面设计             "
$ws.Range("C24").Value = "完成"

# Row 25 (黄俊贤): plan content + completion status
$ws.Range("B25").Value = "初步设计完成网页端的好友添加、创建群主、查询用户界面设计"
$ws.Range("C25").Value = "未完成"

# Row 26 (李达波): plan content + completion status
$ws.Range("B26").Value = "初步设计完成个人信息管理、个人信息修改、找回密码、注册的用例分析和界面设计"
$ws.Range("C26").Value = "完成"

# Row 27 (常永伟 / 收集其他组员所需的图片): completion status only
$ws.Range("C27").Value = "完成"

# Row 28 (冯德志): plan content + completion status
$ws.Range("B28").Value = "初步设计完成网页端的登陆、注册、个人信息管理的用例分析"
$ws.Range("C28").Value = "完成"

# Update the saved view/selection to match the author's last cursor position
$ws.Range("B24").Select()
